$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 22.700661
$ws.Cells.Item(2, 8).Value = 68.10198299999999
$ws.Cells.Item(2, 9).Value = 0.08615268874617349
$ws.Cells.Item(2, 10).Value = 0.08615268874617349
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 13.604331
$ws.Cells.Item(2, 14).Value = 40.812993
$ws.Cells.Item(2, 15).Value = 0.8107276168878804
$ws.Cells.Item(2, 16).Value = 0.8107276168878805
$ws.Cells.Item(2, 17).Value = 308.8273061627909
$ws.Cells.Item(2, 18).Value = 2779.445755465119
$ws.Cells.Item(2, 19).Value = 0.06984636403566855
$ws.Cells.Item(2, 20).Value = 0.06984636403566856

# Row 3
$ws.Cells.Item(3, 7).Value = 22.700661
$ws.Cells.Item(3, 8).Value = 68.10198299999999
$ws.Cells.Item(3, 9).Value = 0.08615268874617349
$ws.Cells.Item(3, 10).Value = 0.08615268874617349
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.061748
$ws.Cells.Item(3, 14).Value = 3.185244
$ws.Cells.Item(3, 15).Value = 0.06327311690486458
$ws.Cells.Item(3, 16).Value = 0.06327311690486459
$ws.Cells.Item(3, 17).Value = 24.102381415428
$ws.Cells.Item(3, 18).Value = 216.921432738852
$ws.Cells.Item(3, 19).Value = 0.005451149146705046
$ws.Cells.Item(3, 20).Value = 0.005451149146705047

# Row 4
$ws.Cells.Item(4, 7).Value = 22.700661
$ws.Cells.Item(4, 8).Value = 68.10198299999999
$ws.Cells.Item(4, 9).Value = 0.08615268874617349
$ws.Cells.Item(4, 10).Value = 0.08615268874617349
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.995771333333333
$ws.Cells.Item(4, 14).Value = 5.987314
$ws.Cells.Item(4, 15).Value = 0.1189346934389115
$ws.Cells.Item(4, 16).Value = 0.1189346934389116
$ws.Cells.Item(4, 17).Value = 45.30532847151799
$ws.Cells.Item(4, 18).Value = 407.7479562436619
$ws.Cells.Item(4, 19).Value = 0.01024654362496411
$ws.Cells.Item(4, 20).Value = 0.01024654362496411

# Row 5
$ws.Cells.Item(5, 7).Value = 22.700661
$ws.Cells.Item(5, 8).Value = 68.10198299999999
$ws.Cells.Item(5, 9).Value = 0.08615268874617349
$ws.Cells.Item(5, 10).Value = 0.08615268874617349
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1185463333333334
$ws.Cells.Item(5, 14).Value = 0.355639
$ws.Cells.Item(5, 15).Value = 0.007064572768343379
$ws.Cells.Item(5, 16).Value = 0.007064572768343379
$ws.Cells.Item(5, 17).Value = 2.691080125793
$ws.Cells.Item(5, 18).Value = 24.219721132137
$ws.Cells.Item(5, 19).Value = 0.0006086319388357804
$ws.Cells.Item(5, 20).Value = 0.0006086319388357804

# Row 6
$ws.Cells.Item(6, 7).Value = 132.5447616666667
$ws.Cells.Item(6, 8).Value = 397.634285
$ws.Cells.Item(6, 9).Value = 0.5030288587986086
$ws.Cells.Item(6, 10).Value = 0.5030288587986087
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 13.604331
$ws.Cells.Item(6, 14).Value = 40.812993
$ws.Cells.Item(6, 15).Value = 0.8107276168878804
$ws.Cells.Item(6, 16).Value = 0.8107276168878805
$ws.Cells.Item(6, 17).Value = 1803.182810029445
$ws.Cells.Item(6, 18).Value = 16228.645290265
$ws.Cells.Item(6, 19).Value = 0.407819387919626
$ws.Cells.Item(6, 20).Value = 0.4078193879196262

# Row 7
$ws.Cells.Item(7, 7).Value = 132.5447616666667
$ws.Cells.Item(7, 8).Value = 397.634285
$ws.Cells.Item(7, 9).Value = 0.5030288587986086
$ws.Cells.Item(7, 10).Value = 0.5030288587986087
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.061748
$ws.Cells.Item(7, 14).Value = 3.185244
$ws.Cells.Item(7, 15).Value = 0.06327311690486458
$ws.Cells.Item(7, 16).Value = 0.06327311690486459
$ws.Cells.Item(7, 17).Value = 140.72913561006
$ws.Cells.Item(7, 18).Value = 1266.56222049054
$ws.Cells.Item(7, 19).Value = 0.03182820378928498
$ws.Cells.Item(7, 20).Value = 0.03182820378928499

# Row 8
$ws.Cells.Item(8, 7).Value = 132.5447616666667
$ws.Cells.Item(8, 8).Value = 397.634285
$ws.Cells.Item(8, 9).Value = 0.5030288587986086
$ws.Cells.Item(8, 10).Value = 0.5030288587986087
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.995771333333333
$ws.Cells.Item(8, 14).Value = 5.987314
$ws.Cells.Item(8, 15).Value = 0.1189346934389115
$ws.Cells.Item(8, 16).Value = 0.1189346934389116
$ws.Cells.Item(8, 17).Value = 264.5290357178322
$ws.Cells.Item(8, 18).Value = 2380.76132146049
$ws.Cells.Item(8, 19).Value = 0.05982758311213802
$ws.Cells.Item(8, 20).Value = 0.05982758311213805

# Row 9
$ws.Cells.Item(9, 7).Value = 132.5447616666667
$ws.Cells.Item(9, 8).Value = 397.634285
$ws.Cells.Item(9, 9).Value = 0.5030288587986086
$ws.Cells.Item(9, 10).Value = 0.5030288587986087
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1185463333333334
$ws.Cells.Item(9, 14).Value = 0.355639
$ws.Cells.Item(9, 15).Value = 0.007064572768343379
$ws.Cells.Item(9, 16).Value = 0.007064572768343379
$ws.Cells.Item(9, 17).Value = 15.71269549812389
$ws.Cells.Item(9, 18).Value = 141.414259483115
$ws.Cells.Item(9, 19).Value = 0.003553683977559497
$ws.Cells.Item(9, 20).Value = 0.003553683977559498

# Row 10
$ws.Cells.Item(10, 7).Value = 41.94534433333333
$ws.Cells.Item(10, 8).Value = 125.836033
$ws.Cells.Item(10, 9).Value = 0.159189381961201
$ws.Cells.Item(10, 10).Value = 0.159189381961201
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 13.604331
$ws.Cells.Item(10, 14).Value = 40.812993
$ws.Cells.Item(10, 15).Value = 0.8107276168878804
$ws.Cells.Item(10, 16).Value = 0.8107276168878805
$ws.Cells.Item(10, 17).Value = 570.638348219641
$ws.Cells.Item(10, 18).Value = 5135.745133976769
$ws.Cells.Item(10, 19).Value = 0.129059228271259
$ws.Cells.Item(10, 20).Value = 0.129059228271259

# Row 11
$ws.Cells.Item(11, 7).Value = 41.94534433333333
$ws.Cells.Item(11, 8).Value = 125.836033
$ws.Cells.Item(11, 9).Value = 0.159189381961201
$ws.Cells.Item(11, 10).Value = 0.159189381961201
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.061748
$ws.Cells.Item(11, 14).Value = 3.185244
$ws.Cells.Item(11, 15).Value = 0.06327311690486458
$ws.Cells.Item(11, 16).Value = 0.06327311690486459
$ws.Cells.Item(11, 17).Value = 44.535385455228
$ws.Cells.Item(11, 18).Value = 400.8184690970519
$ws.Cells.Item(11, 19).Value = 0.01007240837484421
$ws.Cells.Item(11, 20).Value = 0.01007240837484421

# Row 12
$ws.Cells.Item(12, 7).Value = 41.94534433333333
$ws.Cells.Item(12, 8).Value = 125.836033
$ws.Cells.Item(12, 9).Value = 0.159189381961201
$ws.Cells.Item(12, 10).Value = 0.159189381961201
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.995771333333333
$ws.Cells.Item(12, 14).Value = 5.987314
$ws.Cells.Item(12, 15).Value = 0.1189346934389115
$ws.Cells.Item(12, 16).Value = 0.1189346934389116
$ws.Cells.Item(12, 17).Value = 83.71331578726243
$ws.Cells.Item(12, 18).Value = 753.4198420853619
$ws.Cells.Item(12, 19).Value = 0.01893314034228523
$ws.Cells.Item(12, 20).Value = 0.01893314034228524

# Row 13
$ws.Cells.Item(13, 7).Value = 41.94534433333333
$ws.Cells.Item(13, 8).Value = 125.836033
$ws.Cells.Item(13, 9).Value = 0.159189381961201
$ws.Cells.Item(13, 10).Value = 0.159189381961201
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.1185463333333334
$ws.Cells.Item(13, 14).Value = 0.355639
$ws.Cells.Item(13, 15).Value = 0.007064572768343379
$ws.Cells.Item(13, 16).Value = 0.007064572768343379
$ws.Cells.Item(13, 17).Value = 4.972466771120779
$ws.Cells.Item(13, 18).Value = 44.752200940087
$ws.Cells.Item(13, 19).Value = 0.001124604972812513
$ws.Cells.Item(13, 20).Value = 0.001124604972812513

# Row 14
$ws.Cells.Item(14, 7).Value = 66.302588
$ws.Cells.Item(14, 8).Value = 198.907764
$ws.Cells.Item(14, 9).Value = 0.2516290704940168
$ws.Cells.Item(14, 10).Value = 0.2516290704940168
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 13.604331
$ws.Cells.Item(14, 14).Value = 40.812993
$ws.Cells.Item(14, 15).Value = 0.8107276168878804
$ws.Cells.Item(14, 16).Value = 0.8107276168878805
$ws.Cells.Item(14, 17).Value = 902.002353308628
$ws.Cells.Item(14, 18).Value = 8118.021179777651
$ws.Cells.Item(14, 19).Value = 0.2040026366613267
$ws.Cells.Item(14, 20).Value = 0.2040026366613268

# Row 15
$ws.Cells.Item(15, 7).Value = 66.302588
$ws.Cells.Item(15, 8).Value = 198.907764
$ws.Cells.Item(15, 9).Value = 0.2516290704940168
$ws.Cells.Item(15, 10).Value = 0.2516290704940168
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.061748
$ws.Cells.Item(15, 14).Value = 3.185244
$ws.Cells.Item(15, 15).Value = 0.06327311690486458
$ws.Cells.Item(15, 16).Value = 0.06327311690486459
$ws.Cells.Item(15, 17).Value = 70.39664020382399
$ws.Cells.Item(15, 18).Value = 633.569761834416
$ws.Cells.Item(15, 19).Value = 0.01592135559403034
$ws.Cells.Item(15, 20).Value = 0.01592135559403034

# Row 16
$ws.Cells.Item(16, 7).Value = 66.302588
$ws.Cells.Item(16, 8).Value = 198.907764
$ws.Cells.Item(16, 9).Value = 0.2516290704940168
$ws.Cells.Item(16, 10).Value = 0.2516290704940168
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 1.995771333333333
$ws.Cells.Item(16, 14).Value = 5.987314
$ws.Cells.Item(16, 15).Value = 0.1189346934389115
$ws.Cells.Item(16, 16).Value = 0.1189346934389116
$ws.Cells.Item(16, 17).Value = 132.3248044562106
$ws.Cells.Item(16, 18).Value = 1190.923240105896
$ws.Cells.Item(16, 19).Value = 0.02992742635952415
$ws.Cells.Item(16, 20).Value = 0.02992742635952416

# Row 17
$ws.Cells.Item(17, 7).Value = 66.302588
$ws.Cells.Item(17, 8).Value = 198.907764
$ws.Cells.Item(17, 9).Value = 0.2516290704940168
$ws.Cells.Item(17, 10).Value = 0.2516290704940168
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.1185463333333334
$ws.Cells.Item(17, 14).Value = 0.355639
$ws.Cells.Item(17, 15).Value = 0.007064572768343379
$ws.Cells.Item(17, 16).Value = 0.007064572768343379
$ws.Cells.Item(17, 17).Value = 7.859928697910668
$ws.Cells.Item(17, 18).Value = 70.739358281196
$ws.Cells.Item(17, 19).Value = 0.001777651879135588
$ws.Cells.Item(17, 20).Value = 0.001777651879135588
